$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme
Write-Output $cs.Count
for ($i=1; $i -le $cs.Count; $i++) {
    Write-Output "$i : $($cs.Colors.Item($i).RGB)"
}
